$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-32 from serial 45604 to 45605
for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45604) {
        $cell.Value2 = 45605
    }
}
